# Bootcamp-Project-1.pptx — "Added notes to ppt."
#
# 1) Slide 14 speaker notes: the last paragraph ("At-home productivity
#    items, at-home comfort items, at-home entertainment items") was
#    typed as two separate runs with identical formatting; collapse it
#    back into a single run (no text/content change).
# 2) Slide 7 speaker notes: add four new paragraphs after "Lauren"
#    describing the three purchasing-spike categories covered by the
#    following slides (Productivity / Comfort / Entertainment).

$p = $ppt.ActivePresentation

$slide14 = $p.Slides.Item(14)
$notesShape14 = $slide14.NotesPage.Shapes.Item(2)
$notesShape14.TextFrame.TextRange.Text = "Lauren`nConclusions – how to prepare for for an apocalypse (what items to have on hand) and then how to thrive post apocalypse economy`nMention items to barter`n`nNotes from Siara:`nFirst question: Items to include in a Doomsday Bunker:`nFood, Household Cleaning Products, Household Paper Products`nSecond question: Items to buy while others are panicking:`nAt-home productivity items, at-home comfort items, at-home entertainment items"

$slide7 = $p.Slides.Item(7)
$notesShape7 = $slide7.NotesPage.Shapes.Item(2)
$notesShape7.TextFrame.TextRange.Text = "Lauren`nWe noticed purchasing spikes in 3 main categories:`nProductivity`nComfort`nEntertainment"
